# Generate Report for handoff
# Updates the "Latest Handoff Datetime" column (D) for the row corresponding
# to the "8cd77d3d-0189-4bdf-8984-49e473614c01.md" source file on both the
# "zh-cn" and "de-de" worksheets, reflecting a newly generated handoff.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Row 4 on each sheet is the "8cd77d3d-0189-4bdf-8984-49e473614c01.md" entry.
$wsZhCn.Range("D4").Value = "2016-01-17 03:12:37"
$wsDeDe.Range("D4").Value = "2016-01-17 03:12:46"
